$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2026-01-19 01:41:55"

# 1) Rows 2-5 only get their timestamp (column A) refreshed; other columns are unchanged.
$ws.Cells.Item(2, 1).Value = $newTimestamp
$ws.Cells.Item(3, 1).Value = $newTimestamp
$ws.Cells.Item(4, 1).Value = $newTimestamp
$ws.Cells.Item(5, 1).Value = $newTimestamp

# 2) Shift the old row 7 data down to row 8 (content itself is unchanged, only the row moves).
$ws.Cells.Item(8, 1).Value = $newTimestamp
$ws.Cells.Item(8, 2).Value = "初回 PM/PMO(オープン) 一部リモート"
$ws.Cells.Item(8, 3).Value = "システム開発"
$ws.Cells.Item(8, 4).Value = "1,000 ~ 5,000 円 / 固定"
$ws.Cells.Item(8, 5).Value = "期限情報なし"
$ws.Cells.Item(8, 6).Value = "https://www.lancers.jp/work/detail/5473958"
$ws.Cells.Item(8, 6).Style = "Hyperlink"
$ws.Cells.Item(8, 7).Value = 10
$ws.Cells.Item(8, 8).ClearContents()

# 3) Shift the old row 6 data down to row 7 (content itself is unchanged, only the row moves).
$ws.Cells.Item(7, 1).Value = $newTimestamp
$ws.Cells.Item(7, 2).Value = "【急募】プログラム修正依頼!スキルを活かしてみませんか?"
$ws.Cells.Item(7, 3).Value = "システム開発"
$ws.Cells.Item(7, 4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(7, 5).Value = "期限情報なし"
$ws.Cells.Item(7, 6).Value = "https://www.lancers.jp/work/detail/5473840"
$ws.Cells.Item(7, 6).Style = "Hyperlink"
$ws.Cells.Item(7, 7).Value = 13
$ws.Cells.Item(7, 8).ClearContents()

# 4) Write the brand-new listing into row 6.
$ws.Cells.Item(6, 1).Value = $newTimestamp
$ws.Cells.Item(6, 2).Value = "進行管理およびチームディレクションを担当"
$ws.Cells.Item(6, 3).Value = "システム開発"
$ws.Cells.Item(6, 4).Value = "~ 5,000 円 / 固定"
$ws.Cells.Item(6, 5).Value = "期限情報なし"
$ws.Cells.Item(6, 6).Value = "https://www.lancers.jp/work/detail/5418064"
$ws.Cells.Item(6, 6).Style = "Hyperlink"
$ws.Cells.Item(6, 7).Value = 30
$ws.Cells.Item(6, 8).Value = "◇管理"

# 5) Rebuild the hyperlinks collection from scratch for column F so link targets/anchors
#    line up exactly with the rows above (avoids stale hyperlink bindings).
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5473648")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5473858")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5473940")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5468432")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5418064")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5473840")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5473958")

Write-Output "edit complete"
